$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the header in G1 from "Flashcards" to "Flashcard"
$ws.Range("G1").Value = "Flashcard"

# Update the selection to H7 (as recorded in the saved view state)
$ws.Range("H7").Select()
